$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.551.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.806.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.85%  '
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.008'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.98%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4549'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3660'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07127'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8717'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07775'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.21'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.94%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.822.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.280'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.328'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.010'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008572'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.008'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.590.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.72%  '
$ws.Range("E22").Value = '  -3.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.060.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.01%  '
$ws.Range("E24").Value = '  -1.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.983'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.97%  '
$ws.Range("E28").Value = '  -3.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.871'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08689'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.072'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7324'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.85%  '
$ws.Range("E34").Value = '  -1.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.112'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.511'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.073'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01914'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05092'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.863'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.881'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4903'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1568'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.127'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.009'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4593'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.921'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.582'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05997'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.07%  '
